# Trial upload: add a new value below the existing table and move the
# active selection down to reflect the new extent of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: F4 = 777 (extends the used range to A1:F4).
$ws.Range("F4").Value = 777

# Move / update the selection to F5, matching the post-edit cursor position.
$ws.Range("F5").Select()
